$wb = $excel.ActiveWorkbook

# Add the new "PushPull" sheet after the last existing sheet (ReqReply)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "PushPull"

# Row 1
$ws.Range("A1").Value = "tcp: 10000 msgs"

# Row 2
$ws.Range("A2").Value = "Pullers:"
$ws.Range("A2").Font.Bold = $true
$ws.Range("B2").Value = 1
$ws.Range("D2").Value = 2
$ws.Range("F2").Value = 3
$ws.Range("H2").Value = 4

# Row 3 - headers
$ws.Range("A3").Value = "size"
$ws.Range("B3").Value = "msgs/sec"
$ws.Range("C3").Value = "kb/sec"
$ws.Range("D3").Value = "msgs/sec"
$ws.Range("E3").Value = "kb/sec"
$ws.Range("F3").Value = "msgs/sec"
$ws.Range("G3").Value = "kb/sec"
$ws.Range("H3").Value = "msg/sec"
$ws.Range("I3").Value = "kb/sec"

# Match the saved selection / active-cell state on the new sheet
$ws.Range("A4:I5").Select() | Out-Null
